$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Maio")

# Row 12: "Placa de Vídeo GTX650" — move the value from C (Aliexpress) to B,
# and remove the Aliexpress cell in C entirely.
$ws.Range("B12").Value = 390.5
$ws.Range("C12").Clear()

# Currency number format shared by column B's "s=13" cells (e.g. B4, B6..B10)
$currencyFmt = '_-[$R$-416]\ * #,##0.00_-;\-[$R$-416]\ * #,##0.00_-;_-[$R$-416]\ * "-"??_-;_-@_-'

# New row 14: Promofarma
$ws.Range("A14").Value = "Promofarma"
$ws.Range("B14").Value = 118.81
$ws.Range("B14").NumberFormat = $currencyFmt

# New row 15: Loterias
$ws.Range("A15").Value = "Loterias"
$ws.Range("B15").Value = 39
$ws.Range("B15").NumberFormat = $currencyFmt

# New row 16: Steam Games
$ws.Range("A16").Value = "Steam Games"
$ws.Range("B16").Value = 54
$ws.Range("B16").NumberFormat = $currencyFmt

# Update selection to match the diff (active cell moves to A17 after edits)
$ws.Range("A17").Select()
